$wb = $excel.ActiveWorkbook

# --- tradeDetails sheet: trade commencement date changes 01/03/2017 -> 30/03/2017 ---
$ws3 = $wb.Worksheets.Item("tradeDetails")
$ws3.Range("H2").Value = "30/03/2017"

# --- licenseClosure sheet: tradeCategory changes "Veterinary Trades" -> "Flammables" ---
$ws5 = $wb.Worksheets.Item("licenseClosure")
$ws5.Range("C2").Value = "Flammables"
# match the text-cell formatting used elsewhere in the workbook (numFmtId 49 / style index 7)
$ws5.Range("C2").NumberFormat = $ws3.Range("A1").NumberFormat()

# --- update the selected cell on each affected sheet, then leave tradeDetails active ---
$ws3.Activate() | Out-Null
$ws3.Range("D10").Select() | Out-Null

$ws5.Activate() | Out-Null
$ws5.Range("C9").Select() | Out-Null

$ws3.Activate() | Out-Null
